$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.902492
$ws.Range("H2").Value = 3.804984
$ws.Range("M2").Value = 10.402079
$ws.Range("N2").Value = 20.804158
$ws.Range("O2").Value = 0.117441350183963
$ws.Range("P2").Value = 0.08862141909929068
$ws.Range("Q2").Value = 19.789872080868
$ws.Range("R2").Value = 79.159488323472
$ws.Range("S2").Value = 0.117441350183963
$ws.Range("T2").Value = 0.08862141909929068
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.902492
$ws.Range("H3").Value = 3.804984
$ws.Range("O3").Value = 0.6341369869521791
$ws.Range("P3").Value = 0.7177810830557603
$ws.Range("Q3").Value = 106.857506609664
$ws.Range("R3").Value = 641.145039657984
$ws.Range("S3").Value = 0.6341369869521791
$ws.Range("T3").Value = 0.7177810830557603
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.902492
$ws.Range("H4").Value = 3.804984
$ws.Range("M4").Value = 0.2909543333333333
$ws.Range("N4").Value = 0.8728629999999999
$ws.Range("O4").Value = 0.003284926960133785
$ws.Range("P4").Value = 0.003718216220971988
$ws.Range("Q4").Value = 0.553538291532
$ws.Range("R4").Value = 3.321229749192
$ws.Range("S4").Value = 0.003284926960133785
$ws.Range("T4").Value = 0.003718216220971988
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.902492
$ws.Range("H5").Value = 3.804984
$ws.Range("M5").Value = 20.562391
$ws.Range("N5").Value = 41.124782
$ws.Range("O5").Value = 0.232153107282743
$ws.Range("P5").Value = 0.175183083160057
$ws.Range("Q5").Value = 39.119784378372
$ws.Range("R5").Value = 156.479137513488
$ws.Range("S5").Value = 0.232153107282743
$ws.Range("T5").Value = 0.175183083160057
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.902492
$ws.Range("H6").Value = 3.804984
$ws.Range("M6").Value = 0.5741476666666667
$ws.Range("N6").Value = 1.722443
$ws.Range("O6").Value = 0.006482230828885768
$ws.Range("P6").Value = 0.007337251667557973
$ws.Range("Q6").Value = 1.092311342652
$ws.Range("R6").Value = 6.553868055912001
$ws.Range("S6").Value = 0.006482230828885768
$ws.Range("T6").Value = 0.007337251667557973
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.902492
$ws.Range("H7").Value = 3.804984
$ws.Range("M7").Value = 0.5758453333333333
$ws.Range("N7").Value = 1.727536
$ws.Range("O7").Value = 0.006501397792095299
$ws.Range("P7").Value = 0.00735894679636216
$ws.Range("Q7").Value = 1.095541139904
$ws.Range("R7").Value = 6.573246839424
$ws.Range("S7").Value = 0.006501397792095299
$ws.Range("T7").Value = 0.00735894679636216
